$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-cell formatting from the last existing row (A238) down
# through the new rows (A239:A244) so the new date cells pick up the
# same style (centered, bordered, custom date format) as the rest of
# column A.
$ws.Range("A238").Copy()
$ws.Range("A239:A244").PasteSpecial(-4122)

# New daily data rows covering 2021-04-27 through 2021-05-02
# (serials 44313-44318).
$data = @(
    @(44313, 1, 2, 50.8646998982706),
    @(44314, 0, 2, 50.8646998982706),
    @(44315, 0, 2, 50.8646998982706),
    @(44316, 1, 3, 76.2970498474059),
    @(44317, 3, 6, 152.5940996948118),
    @(44318, 0, 5, 127.1617497456765)
)

$r = 239
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
